# Update cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51.
# Numeric-looking price strings are prefixed with a literal apostrophe so Excel
# stores them as text (preserving formats like "41.50" / "0.05200") instead of
# silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.759.81"
$ws.Range("E2").Value = "  -2.82%  "
# Row 3
$ws.Range("D3").Value = "1.746.09"
$ws.Range("E3").Value = "  -5.26%  "
# Row 4
$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  -0.45%  "
# Row 5
$ws.Range("D5").Value = "'238.17"
$ws.Range("E5").Value = "  -9.20%  "
# Row 6
$ws.Range("D6").Value = "'0.9974"
$ws.Range("E6").Value = "  -0.35%  "
# Row 7
$ws.Range("D7").Value = "'0.4949"
$ws.Range("E7").Value = "  -6.86%  "
# Row 8
$ws.Range("D8").Value = "'41.64"
$ws.Range("E8").Value = "  -7.32%  "
# Row 9
$ws.Range("D9").Value = "'0.2654"
$ws.Range("E9").Value = "  -13.64%  "
# Row 10
$ws.Range("D10").Value = "'0.06107"
$ws.Range("E10").Value = "  -11.52%  "
# Row 11
$ws.Range("D11").Value = "1.743.05"
$ws.Range("E11").Value = "  -5.37%  "
# Row 12
$ws.Range("D12").Value = "'0.06899"
$ws.Range("E12").Value = "  -11.57%  "
# Row 13
$ws.Range("D13").Value = "'15.31"
$ws.Range("E13").Value = "  -16.91%  "
# Row 14
$ws.Range("D14").Value = "'4.476"
$ws.Range("E14").Value = "  -10.78%  "
# Row 15
$ws.Range("D15").Value = "'76.80"
$ws.Range("E15").Value = "  -14.37%  "
# Row 16
$ws.Range("D16").Value = "'0.5844"
$ws.Range("E16").Value = "  -22.74%  "
# Row 17
$ws.Range("D17").Value = "'0.9955"
$ws.Range("E17").Value = "  -0.53%  "
# Row 18
$ws.Range("D18").Value = "'0.9981"
$ws.Range("E18").Value = "  -0.28%  "
# Row 19
$ws.Range("D19").Value = "25.801.23"
$ws.Range("E19").Value = "  -2.74%  "
# Row 20
$ws.Range("D20").Value = "'11.59"
$ws.Range("E20").Value = "  -17.28%  "
# Row 21
$ws.Range("D21").Value = "'0.000006682"
$ws.Range("E21").Value = "  -15.93%  "
# Row 22
$ws.Range("D22").Value = "1.961.75"
$ws.Range("E22").Value = "  -5.49%  "
# Row 23
$ws.Range("D23").Value = "'4.048"
# Row 24
$ws.Range("D24").Value = "'7.984"
$ws.Range("E24").Value = "  -14.35%  "
# Row 25
$ws.Range("D25").Value = "'5.070"
$ws.Range("E25").Value = "  -15.32%  "
# Row 26
$ws.Range("D26").Value = "'137.60"
$ws.Range("E26").Value = "  -3.32%  "
# Row 27
$ws.Range("D27").Value = "'1.523"
$ws.Range("E27").Value = "  -9.82%  "
# Row 28
$ws.Range("D28").Value = "'1.839"
$ws.Range("E28").Value = "  -16.17%  "
# Row 29
$ws.Range("D29").Value = "'14.84"
$ws.Range("E29").Value = "  -12.62%  "
# Row 30
$ws.Range("D30").Value = "'102.06"
$ws.Range("E30").Value = "  -8.11%  "
# Row 31
$ws.Range("D31").Value = "'3.745"
$ws.Range("E31").Value = "  -12.39%  "
# Row 32
$ws.Range("D32").Value = "'0.08020"
$ws.Range("E32").Value = "  -8.86%  "
# Row 33
$ws.Range("D33").Value = "'3.476"
$ws.Range("E33").Value = "  -15.06%  "
# Row 34
$ws.Range("D34").Value = "'0.04444"
$ws.Range("E34").Value = "  -7.83%  "
# Row 35
$ws.Range("D35").Value = "'0.9954"
$ws.Range("E35").Value = "  -0.45%  "
# Row 36
$ws.Range("D36").Value = "'2.623"
$ws.Range("E36").Value = "  -10.51%  "
# Row 37
$ws.Range("D37").Value = "'0.9793"
$ws.Range("E37").Value = "  -13.62%  "
# Row 38
$ws.Range("D38").Value = "'0.5989"
$ws.Range("E38").Value = "  -17.91%  "
# Row 39
$ws.Range("D39").Value = "'2.650"
$ws.Range("E39").Value = "  -14.74%  "
# Row 40
$ws.Range("D40").Value = "'104.98"
$ws.Range("E40").Value = "  -2.97%  "
# Row 41
$ws.Range("D41").Value = "'1.926"
$ws.Range("E41").Value = "  -16.68%  "
# Row 42
$ws.Range("D42").Value = "'0.9967"
$ws.Range("E42").Value = "  -0.38%  "
# Row 43
$ws.Range("D43").Value = "'0.01516"
$ws.Range("E43").Value = "  -11.69%  "
# Row 44
$ws.Range("D44").Value = "'5.149"
$ws.Range("E44").Value = "  -12.23%  "
# Row 45
$ws.Range("D45").Value = "'0.3794"
$ws.Range("E45").Value = "  -20.98%  "
# Row 46
$ws.Range("D46").Value = "'0.7252"
$ws.Range("E46").Value = "  -19.63%  "
# Row 47
$ws.Range("D47").Value = "'0.1110"
$ws.Range("E47").Value = "  -10.69%  "
# Row 48
$ws.Range("D48").Value = "'0.05196"
$ws.Range("E48").Value = "  -10.36%  "
# Row 49
$ws.Range("D49").Value = "'30.08"
$ws.Range("E49").Value = "  -13.71%  "
# Row 50
$ws.Range("D50").Value = "'5.907"
$ws.Range("E50").Value = "  -21.14%  "
# Row 51
$ws.Range("D51").Value = "'52.24"
$ws.Range("E51").Value = "  -13.31%  "
